$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlEdgeTop = 8
$xlLineStyleNone = -4142
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Insert a new blank row above the old "Total:" row (old row 17), pushing
#    it (and everything below) down by one. This is the "added a date row"
#    from the commit message.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Insert()

# ---------------------------------------------------------------------------
# 2) Fix up the formatting of the new row 17 while row 16 still carries its
#    original ("closing row of the table") formatting, so we can harvest it.
#    Row 17 becomes the new row directly above the totals row, so it takes on
#    a look that's a blend of a normal data row and the old totals-row border.
# ---------------------------------------------------------------------------
$ws.Range("B9").Copy()
$ws.Range("B17").PasteSpecial($xlPasteFormats)

$ws.Range("E9").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)

$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial($xlPasteFormats)

$ws.Range("F16").Copy()
$ws.Range("F17").PasteSpecial($xlPasteFormats)
$ws.Range("J17").PasteSpecial($xlPasteFormats)

$ws.Range("H16").Copy()
$ws.Range("H17").PasteSpecial($xlPasteFormats)

$ws.Range("L16").Copy()
$ws.Range("L17").PasteSpecial($xlPasteFormats)

$ws.Range("G18").Copy()
$ws.Range("G17").PasteSpecial($xlPasteFormats)

# Date columns E17/I17/K17: start from a normal full-box date cell, then
# drop just the top edge so it reads as "left+right+bottom" only.
$ws.Range("E9").Copy()
$ws.Range("E17").PasteSpecial($xlPasteFormats)
$ws.Range("I17").PasteSpecial($xlPasteFormats)
$ws.Range("K17").PasteSpecial($xlPasteFormats)
$ws.Range("E17").Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone
$ws.Range("I17").Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone
$ws.Range("K17").Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone

# ---------------------------------------------------------------------------
# 3) Now "open up" row 16 so it reads as a normal middle row of the table
#    instead of the table's closing/bottom-bordered row (that look moved to
#    the new row 17 above).
# ---------------------------------------------------------------------------
$ws.Range("B9").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("D16").PasteSpecial($xlPasteFormats)
$ws.Range("F16").PasteSpecial($xlPasteFormats)
$ws.Range("H16").PasteSpecial($xlPasteFormats)

$ws.Range("E9").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("E16").PasteSpecial($xlPasteFormats)
$ws.Range("G16").PasteSpecial($xlPasteFormats)
$ws.Range("I16").PasteSpecial($xlPasteFormats)
$ws.Range("K16").PasteSpecial($xlPasteFormats)

$ws.Range("J10").Copy()
$ws.Range("J16").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Fill in the new date/hours entry added on row 15 (Andrew's column).
# ---------------------------------------------------------------------------
$ws.Range("E15").Value = 44478
$ws.Range("F15").Value = 4

# ---------------------------------------------------------------------------
# 5) Extend the SUM ranges on the totals row (now row 18) to include the
#    newly inserted row 17.
# ---------------------------------------------------------------------------
$ws.Range("D18").Formula = "=SUM(D9:D17)"
$ws.Range("F18").Formula = "=SUM(F9:F17)"
$ws.Range("H18").Formula = "=SUM(H9:H17)"
$ws.Range("J18").Formula = "=SUM(J9:J17)"

# ---------------------------------------------------------------------------
# 6) View tweaks captured in the diff (zoom + new selection).
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 120
$ws.Range("F21").Select()
